$d = $word.ActiveDocument

# Paragraph 1: 1G port pricing - two dollar amounts to convert ($0 and $250)
$d.Content.Find.Execute(
    "A 1G port is `$0/year. More than one 1G port requires board approval; if the board grants an exception, the participant is charged the 10G fee. That is, a 2x1G LAG would be `$250/year.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A 1G port is 0 USD/year. More than one 1G port requires board approval; if the board grants an exception, the participant is charged the 10G fee. That is, a 2x1G LAG would be 250 USD/year.",
    2) | Out-Null

# Paragraph 2: 10G port pricing - two dollar amounts to convert ($250 and $1,000)
$d.Content.Find.Execute(
    "The first 10G port is `$250/year. Additional 10G ports are `$1,000/year. More than four ports requires board approval.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The first 10G port is 250 USD/year. Additional 10G ports are 1,000 USD/year. More than four ports requires board approval.",
    2) | Out-Null

# Paragraph 3: 100G port pricing - one dollar amount to convert ($3,000)
$d.Content.Find.Execute(
    "A 100G port is `$3,000/year. We currently have no limit for 100G ports besides available ports and reasonable technical need.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A 100G port is 3,000 USD/year. We currently have no limit for 100G ports besides available ports and reasonable technical need.",
    2) | Out-Null
